$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "phone_no" column header (F1) and the shared-string table gains it.
$ws.Range("F1").Value = "phone_no"

# First few rows hold literal phone numbers (not part of the incrementing series).
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 67
$ws.Range("F4").Value = 7387256694
$ws.Range("F5").Value = 1

# Remaining rows were filled down as an incrementing series (F5+1, F6+1, ...).
$ws.Range("F6:F69").Formula = "=F5+1"
$ws.Range("F70").Formula = "=F69+1"

# New column sized to fit its content, like the other data columns.
$ws.Columns.Item(6).ColumnWidth = 10.140625

# Selection moved to F3 after the edit, matching the saved view state.
$ws.Range("F3").Select()

Write-Host "done"
